$wb = $excel.ActiveWorkbook

# --- Reorder tabs: "review_info" first, "hotel_info" second -------------
$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")
$hotelSheet.Move($null, $reviewSheet)

# --- Insert a new "State" column into hotel_info, between Hotel_Name and City ---
$ws = $wb.Worksheets.Item("hotel_info")
$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "State"
$ws.Range("C2").Value = "Louisiana"
